# Applies the 2026-01-13 05:18:56 scrape refresh to the 3 'Linea 141' schedule sheets.
# For each sheet: bump the header timestamp/row-count, append the chronologically-last
# new trips at the bottom, then insert the remaining new trips at their correct
# chronological position (processed bottom-to-top so row numbers stay stable).

$wb = $excel.ActiveWorkbook

# ===================== Sheet 1: sheet1 =====================
$ws = $wb.Worksheets.Item(1)
$ws.Range('A2').Value = 'Última actualización: 05:18:56'
$ws.Range('A3').Value = 'Total filas: 37'

# Append new trips after the current last data row (31)
$ws.Cells.Item(32, 1).Value = '05:18:56'
$ws.Cells.Item(32, 2).Value = '06:58'
$ws.Cells.Item(32, 3).Value = '10_OLMOS'
$ws.Cells.Item(32, 4).Value = 100
$ws.Cells.Item(32, 5).Value = 'LP1912'

$ws.Cells.Item(33, 1).Value = '05:18:56'
$ws.Cells.Item(33, 2).Value = '06:59'
$ws.Cells.Item(33, 3).Value = '14_ABASTO'
$ws.Cells.Item(33, 4).Value = 101
$ws.Cells.Item(33, 5).Value = 'LP1912'

$ws.Cells.Item(34, 1).Value = '05:18:56'
$ws.Cells.Item(34, 2).Value = '07:04'
$ws.Cells.Item(34, 3).Value = '15_ABASTO'
$ws.Cells.Item(34, 4).Value = 106
$ws.Cells.Item(34, 5).Value = 'LP1912'

$ws.Cells.Item(35, 1).Value = '05:18:56'
$ws.Cells.Item(35, 2).Value = '07:06'
$ws.Cells.Item(35, 3).Value = '225_GOMEZ'
$ws.Cells.Item(35, 4).Value = 108
$ws.Cells.Item(35, 5).Value = 'LP1912'

$ws.Cells.Item(36, 1).Value = '05:18:56'
$ws.Cells.Item(36, 2).Value = '07:11'
$ws.Cells.Item(36, 3).Value = '215A_EL PATO'
$ws.Cells.Item(36, 4).Value = 113
$ws.Cells.Item(36, 5).Value = 'LP1912'

$ws.Cells.Item(37, 1).Value = '05:18:56'
$ws.Cells.Item(37, 2).Value = '07:15'
$ws.Cells.Item(37, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(37, 4).Value = 117
$ws.Cells.Item(37, 5).Value = 'LP1912'

# Insert remaining new trips at their chronological position
# (processed from the bottom-most original row upward so earlier
#  original row numbers below remain valid row references)
$ws.Rows(31).Insert()
$ws.Cells.Item(31, 1).Value = '05:18:56'
$ws.Cells.Item(31, 2).Value = '06:46'
$ws.Cells.Item(31, 3).Value = '215C_EL PATO'
$ws.Cells.Item(31, 4).Value = 88
$ws.Cells.Item(31, 5).Value = 'LP1912'

$ws.Rows(30).Insert()
$ws.Cells.Item(30, 1).Value = '05:18:56'
$ws.Cells.Item(30, 2).Value = '06:43'
$ws.Cells.Item(30, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(30, 4).Value = 85
$ws.Cells.Item(30, 5).Value = 'LP1912'

$ws.Rows(26).Insert()
$ws.Cells.Item(26, 1).Value = '05:18:56'
$ws.Cells.Item(26, 2).Value = '06:26'
$ws.Cells.Item(26, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(26, 4).Value = 68
$ws.Cells.Item(26, 5).Value = 'LP1912'

$ws.Rows(25).Insert()
$ws.Cells.Item(25, 1).Value = '05:18:56'
$ws.Cells.Item(25, 2).Value = '06:20'
$ws.Cells.Item(25, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(25, 4).Value = 62
$ws.Cells.Item(25, 5).Value = 'LP1912'

$ws.Rows(13).Insert()
$ws.Cells.Item(13, 1).Value = '05:18:56'
$ws.Cells.Item(13, 2).Value = '05:20'
$ws.Cells.Item(13, 3).Value = '14_ABASTO'
$ws.Cells.Item(13, 4).Value = 2
$ws.Cells.Item(13, 5).Value = 'LP1912'

# ===================== Sheet 2: sheet2 =====================
$ws = $wb.Worksheets.Item(2)
$ws.Range('A2').Value = 'Última actualización: 05:18:56'
$ws.Range('A3').Value = 'Total filas: 9'

# Append new trips after the current last data row (12)
$ws.Cells.Item(13, 1).Value = '05:18:56'
$ws.Cells.Item(13, 2).Value = '07:11'
$ws.Cells.Item(13, 3).Value = '215A_EL PATO'
$ws.Cells.Item(13, 4).Value = 113
$ws.Cells.Item(13, 5).Value = 'LP1912'

# Insert remaining new trips at their chronological position
# (processed from the bottom-most original row upward so earlier
#  original row numbers below remain valid row references)
$ws.Rows(12).Insert()
$ws.Cells.Item(12, 1).Value = '05:18:56'
$ws.Cells.Item(12, 2).Value = '06:46'
$ws.Cells.Item(12, 3).Value = '215C_EL PATO'
$ws.Cells.Item(12, 4).Value = 88
$ws.Cells.Item(12, 5).Value = 'LP1912'

# ===================== Sheet 3: sheet3 =====================
$ws = $wb.Worksheets.Item(3)
$ws.Range('A2').Value = 'Última actualización: 05:18:56'
$ws.Range('A3').Value = 'Total filas: 7'

# Append new trips after the current last data row (10)
$ws.Cells.Item(11, 1).Value = '05:18:56'
$ws.Cells.Item(11, 2).Value = '06:59'
$ws.Cells.Item(11, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(11, 4).Value = 101
$ws.Cells.Item(11, 5).Value = 'L6173'

# Insert remaining new trips at their chronological position
# (processed from the bottom-most original row upward so earlier
#  original row numbers below remain valid row references)
$ws.Rows(10).Insert()
$ws.Cells.Item(10, 1).Value = '05:18:56'
$ws.Cells.Item(10, 2).Value = '06:32'
$ws.Cells.Item(10, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(10, 4).Value = 74
$ws.Cells.Item(10, 5).Value = 'L6203'

